# Uploaded doorlift door brackets
#
# Inserts two new "Chamber" parts ("Doorlift Door Bracket A_1x" and
# "Doorlift Door Bracket B_1x") into the printed-parts list on Sheet1,
# right after the existing "Door Window Frame Outer_1x" row (old row 40),
# which pushes every row below (old rows 41-68) down by two. Also flips
# the Support flag for the "AB Toolhead Carrier_1x" row from "No" to "BPO".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert two blank rows right after row 40 (shifts 41.. down to 43..) ---
$insertRange = $ws.Range("A41:H42")
$insertRange.Insert()

# --- Row 41: Doorlift Door Bracket A_1x ---
$ws.Range("A41").Value2 = "Chamber"
$ws.Range("B41").Value2 = "Doorlift Door Bracket A_1x"
$ws.Range("C41").Value2 = 0.083333333333333329
$ws.Range("D41").Value2 = 80
$ws.Range("E41").Value2 = 1
$ws.Range("F41").Formula = "=E41*C41"
$ws.Range("G41").Formula = "=E41*D41"
$ws.Range("H41").Value2 = "No"

# --- Row 42: Doorlift Door Bracket B_1x ---
$ws.Range("A42").Value2 = "Chamber"
$ws.Range("B42").Value2 = "Doorlift Door Bracket B_1x"
$ws.Range("C42").Value2 = 0.083333333333333329
$ws.Range("D42").Value2 = 80
$ws.Range("E42").Value2 = 1
$ws.Range("F42").Formula = "=E42*C42"
$ws.Range("G42").Formula = "=E42*D42"
$ws.Range("H42").Value2 = "No"

# --- Unrelated single-cell edit: AB Toolhead Carrier_1x Support -> BPO ---
$ws.Range("H11").Value2 = "BPO"

# --- Keep the hidden AutoFilter defined name in sync with the new range ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$4:`$H`$65"
    }
}

Write-Output "Inserted Doorlift Door Bracket A/B rows; updated H11 and filter range."
